$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 603.6923
$ws.Range("I33").Value = 673.4545000000001
$ws.Range("J33").Value = 220
$ws.Range("K33").Value = 673.4545000000001
$ws.Range("L33").Value = 220
$ws.Range("M33").Value = -444.4545000000001
$ws.Range("N33").Value = -678
$ws.Range("H51").Value = 1854.1428
$ws.Range("I51").Value = 1913.1666
$ws.Range("K51").Value = 1913.1666
$ws.Range("M51").Value = -1429.1666
$ws.Range("H97").Value = 801.9
$ws.Range("I97").Value = 699
$ws.Range("J97").Value = 813.3333
$ws.Range("K97").Value = 2097
$ws.Range("L97").Value = 2439.9999
$ws.Range("M97").Value = -1601
$ws.Range("N97").Value = -3431.9999
$ws.Range("H113").Value = 2906.0715
$ws.Range("I113").Value = 2807
$ws.Range("J113").Value = 2961.111
$ws.Range("K113").Value = 2807
$ws.Range("L113").Value = 2961.111
$ws.Range("M113").Value = 447
$ws.Range("N113").Value = -9469.111000000001
$ws.Range("H132").Value = 3138.8538
$ws.Range("I132").Value = 2875.6943
$ws.Range("J132").Value = 5033.6
$ws.Range("K132").Value = 8627.082900000001
$ws.Range("L132").Value = 15100.8
$ws.Range("M132").Value = -6097.082900000001
$ws.Range("N132").Value = -20160.8

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2022.7273
$ws.Range("I97").Value = 1578.5714
$ws.Range("J97").Value = 2800
$ws.Range("K97").Value = 1578.5714
$ws.Range("L97").Value = 2800
$ws.Range("M97").Value = -1082.5714
$ws.Range("N97").Value = -3792
$ws.Range("H122").Value = 42309.76
$ws.Range("I122").Value = 73246
$ws.Range("J122").Value = 2936.3635
$ws.Range("K122").Value = 219738
$ws.Range("L122").Value = 8809.0905
$ws.Range("M122").Value = -217288
$ws.Range("N122").Value = -13709.0905
$ws.Range("H132").Value = 3461.9756
$ws.Range("I132").Value = 2169.7036
$ws.Range("J132").Value = 5954.2144
$ws.Range("K132").Value = 6509.110799999999
$ws.Range("L132").Value = 17862.6432
$ws.Range("M132").Value = -3979.110799999999
$ws.Range("N132").Value = -22922.6432

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 52633620
$ws.Range("I20").Value = 1873
$ws.Range("J20").Value = 142859470
$ws.Range("K20").Value = 1873
$ws.Range("L20").Value = 142859470
$ws.Range("M20").Value = -1626
$ws.Range("N20").Value = -142859964
$ws.Range("H134").Value = 3544.8438
$ws.Range("I134").Value = 3574.318
$ws.Range("K134").Value = 10722.954
$ws.Range("M134").Value = -8187.954000000002

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4899.6313
$ws.Range("I31").Value = 1386.9166
$ws.Range("J31").Value = 10921.429
$ws.Range("K31").Value = 1386.9166
$ws.Range("L31").Value = 10921.429
$ws.Range("M31").Value = -1091.9166
$ws.Range("N31").Value = -11511.429
$ws.Range("H34").Value = 4899.6313
$ws.Range("I34").Value = 1386.9166
$ws.Range("J34").Value = 10921.429
$ws.Range("K34").Value = 1386.9166
$ws.Range("L34").Value = 10921.429
$ws.Range("M34").Value = -1184.9166
$ws.Range("N34").Value = -11325.429
$ws.Range("H43").Value = 40000
$ws.Range("J43").Value = 40000
$ws.Range("L43").Value = 40000
$ws.Range("N43").Value = -40368
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -11316
$ws.Range("H99").Value = 2052.6956
$ws.Range("I99").Value = 1356
$ws.Range("J99").Value = 2119.0476
$ws.Range("K99").Value = 1356
$ws.Range("L99").Value = 2119.0476
$ws.Range("M99").Value = 142
$ws.Range("N99").Value = -5115.0476
$ws.Range("H101").Value = 40000
$ws.Range("J101").Value = 40000
$ws.Range("L101").Value = 40000
$ws.Range("N101").Value = -46490
$ws.Range("H107").Value = 2404612.8
$ws.Range("I107").Value = 4464699
$ws.Range("K107").Value = 4464699
$ws.Range("M107").Value = -4462779
$ws.Range("H126").Value = 2052.6956
$ws.Range("I126").Value = 1356
$ws.Range("J126").Value = 2119.0476
$ws.Range("K126").Value = 4068
$ws.Range("L126").Value = 6357.1428
$ws.Range("M126").Value = -1598
$ws.Range("N126").Value = -11297.1428
$ws.Range("H134").Value = 6684.1816
$ws.Range("I134").Value = 12037.4
$ws.Range("J134").Value = 2223.1667
$ws.Range("K134").Value = 36112.2
$ws.Range("L134").Value = 6669.500100000001
$ws.Range("M134").Value = -33577.2
$ws.Range("N134").Value = -11739.5001

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1721.9
$ws.Range("I5").Value = 955
$ws.Range("J5").Value = 2233.1667
$ws.Range("K5").Value = 2865
$ws.Range("L5").Value = 6699.500100000001
$ws.Range("M5").Value = -2753
$ws.Range("N5").Value = -6923.500100000001
$ws.Range("H14").Value = 111111440
$ws.Range("I14").Value = 111111440
$ws.Range("K14").Value = 333334320
$ws.Range("M14").Value = -333334147
$ws.Range("H109").Value = 2011.9131
$ws.Range("I109").Value = 758.3333
$ws.Range("J109").Value = 2817.7856
$ws.Range("K109").Value = 2274.9999
$ws.Range("L109").Value = 8453.356800000001
$ws.Range("M109").Value = -1234.9999
$ws.Range("N109").Value = -10533.3568
$ws.Range("H121").Value = 1202.3928
$ws.Range("I121").Value = 465
$ws.Range("J121").Value = 1259.1154
$ws.Range("K121").Value = 1395
$ws.Range("L121").Value = 3777.3462
$ws.Range("M121").Value = -85
$ws.Range("N121").Value = -6397.3462
$ws.Range("H129").Value = 1733.4231
$ws.Range("I129").Value = 826
$ws.Range("J129").Value = 1949.4762
$ws.Range("K129").Value = 2478
$ws.Range("L129").Value = 5848.4286
$ws.Range("M129").Value = 2522
$ws.Range("N129").Value = -15848.4286
$ws.Range("H131").Value = 4208.811
$ws.Range("J131").Value = 4816.4688
$ws.Range("L131").Value = 14449.4064
$ws.Range("N131").Value = -24529.4064
$ws.Range("H134").Value = 6678.3076
$ws.Range("I134").Value = 4582.25
$ws.Range("J134").Value = 7609.8887
$ws.Range("K134").Value = 13746.75
$ws.Range("L134").Value = 22829.6661
$ws.Range("M134").Value = -8676.75
$ws.Range("N134").Value = -32969.6661
$ws.Range("H135").Value = 1721.9
$ws.Range("I135").Value = 955
$ws.Range("J135").Value = 2233.1667
$ws.Range("K135").Value = 8595
$ws.Range("L135").Value = 20098.5003
$ws.Range("M135").Value = -6060
$ws.Range("N135").Value = -25168.5003
$ws.Range("H136").Value = 3540
$ws.Range("I136").Value = 1800
$ws.Range("J136").Value = 3975
$ws.Range("K136").Value = 5400
$ws.Range("L136").Value = 11925
$ws.Range("M136").Value = -300
$ws.Range("N136").Value = -22125
$ws.Range("H137").Value = 49090.75
$ws.Range("J137").Value = 59466.61
$ws.Range("L137").Value = 178399.83
$ws.Range("N137").Value = -188599.83
$ws.Range("H139").Value = 3083.4285
$ws.Range("I139").Value = 2358.2
$ws.Range("J139").Value = 3742.7273
$ws.Range("K139").Value = 7074.599999999999
$ws.Range("L139").Value = 11228.1819
$ws.Range("M139").Value = -1934.599999999999
$ws.Range("N139").Value = -21508.1819
$ws.Range("H140").Value = 1876.6666
$ws.Range("I140").Value = 1592.7273
$ws.Range("K140").Value = 4778.1819
$ws.Range("M140").Value = 401.8181000000004

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1177162.4
$ws.Range("J80").Value = 69133.164
$ws.Range("L80").Value = 69133.164
$ws.Range("N80").Value = -71129.164
$ws.Range("H83").Value = 1177162.4
$ws.Range("J83").Value = 69133.164
$ws.Range("L83").Value = 345665.82
$ws.Range("N83").Value = -355649.82
$ws.Range("H102").Value = 2447
$ws.Range("I102").Value = 2562.4
$ws.Range("J102").Value = 2374.875
$ws.Range("K102").Value = 2562.4
$ws.Range("L102").Value = 2374.875
$ws.Range("M102").Value = -940.4000000000001
$ws.Range("N102").Value = -5618.875
$ws.Range("H132").Value = 7064.905
$ws.Range("I132").Value = 7549.1763
$ws.Range("J132").Value = 5006.75
$ws.Range("K132").Value = 22647.5289
$ws.Range("L132").Value = 15020.25
$ws.Range("M132").Value = -20117.5289
$ws.Range("N132").Value = -20080.25

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1813.8572
$ws.Range("I46").Value = 1799.25
$ws.Range("J46").Value = 1833.3334
$ws.Range("K46").Value = 1799.25
$ws.Range("L46").Value = 1833.3334
$ws.Range("M46").Value = -1611.25
$ws.Range("N46").Value = -2209.3334
$ws.Range("H61").Value = 2356.95
$ws.Range("I61").Value = 1648.9333
$ws.Range("K61").Value = 1648.9333
$ws.Range("M61").Value = -1446.9333
$ws.Range("H113").Value = 2356.95
$ws.Range("I113").Value = 1648.9333
$ws.Range("K113").Value = 1648.9333
$ws.Range("M113").Value = 521.0667000000001
$ws.Range("H132").Value = 4283.5557
$ws.Range("I132").Value = 3689.3333
$ws.Range("J132").Value = 4877.778
$ws.Range("K132").Value = 11067.9999
$ws.Range("L132").Value = 14633.334
$ws.Range("M132").Value = -8537.999899999999
$ws.Range("N132").Value = -19693.334

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 32000
$ws.Range("J75").Value = 32000
$ws.Range("L75").Value = 32000
$ws.Range("N75").Value = -33872
$ws.Range("H78").Value = 32000
$ws.Range("J78").Value = 32000
$ws.Range("L78").Value = 96000
$ws.Range("N78").Value = -105360
$ws.Range("H122").Value = 2691.4736
$ws.Range("I122").Value = 2702.375
$ws.Range("J122").Value = 2633.3333
$ws.Range("K122").Value = 8107.125
$ws.Range("L122").Value = 7899.999899999999
$ws.Range("M122").Value = -5657.125
$ws.Range("N122").Value = -12799.9999
$ws.Range("H132").Value = 6176270.5
$ws.Range("I132").Value = 3789.2666
$ws.Range("J132").Value = 13891872
$ws.Range("K132").Value = 11367.7998
$ws.Range("L132").Value = 41675616
$ws.Range("M132").Value = -8837.799800000001
$ws.Range("N132").Value = -41680676
